$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 1345.0526
$ws.Range("J29").Value = 1792.7142
$ws.Range("L29").Value = 5378.142599999999
$ws.Range("N29").Value = -5940.142599999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 2503.1
$ws.Range("I38").Value = 125.375
$ws.Range("J38").Value = 4088.25
$ws.Range("K38").Value = 376.125
$ws.Range("L38").Value = 12264.75
$ws.Range("M38").Value = -4.125
$ws.Range("N38").Value = -13008.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 2470.4285
$ws.Range("I58").Value = 123.25
$ws.Range("J58").Value = 5600
$ws.Range("K58").Value = 369.75
$ws.Range("L58").Value = 16800
$ws.Range("M58").Value = -219.75
$ws.Range("N58").Value = -17100

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H82").Value = 9575.1875
$ws.Range("I82").Value = 1720.3
$ws.Range("J82").Value = 22666.666
$ws.Range("K82").Value = 5160.9
$ws.Range("L82").Value = 67999.99800000001
$ws.Range("M82").Value = -4754.9
$ws.Range("N82").Value = -68811.99800000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H85").Value = 9575.1875
$ws.Range("I85").Value = 1720.3
$ws.Range("J85").Value = 22666.666
$ws.Range("K85").Value = 5160.9
$ws.Range("L85").Value = 67999.99800000001
$ws.Range("M85").Value = -3756.9
$ws.Range("N85").Value = -70807.99800000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 36328.1
$ws.Range("J87").Value = 36328.1
$ws.Range("L87").Value = 36328.1
$ws.Range("N87").Value = -38824.1

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H90").Value = 36328.1
$ws.Range("J90").Value = 36328.1
$ws.Range("L90").Value = 108984.3
$ws.Range("N90").Value = -121464.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2966.125
$ws.Range("I2").Value = 1247.6
$ws.Range("J2").Value = 5830.3335
$ws.Range("K2").Value = 1247.6
$ws.Range("L2").Value = 5830.3335
$ws.Range("M2").Value = -1134.6
$ws.Range("N2").Value = -6056.3335

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2787.2354
$ws.Range("I45").Value = 2637.1538
$ws.Range("J45").Value = 3275
$ws.Range("K45").Value = 2637.1538
$ws.Range("L45").Value = 3275
$ws.Range("M45").Value = -2260.1538
$ws.Range("N45").Value = -4029

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H101").Value = 32998
$ws.Range("J101").Value = 32998
$ws.Range("L101").Value = 32998
$ws.Range("N101").Value = -39488

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 2966.125
$ws.Range("I116").Value = 1247.6
$ws.Range("J116").Value = 5830.3335
$ws.Range("K116").Value = 1247.6
$ws.Range("L116").Value = 5830.3335
$ws.Range("M116").Value = 1046.4
$ws.Range("N116").Value = -10418.3335

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2875.12
$ws.Range("I3").Value = 1212.8125
$ws.Range("J3").Value = 5830.3335
$ws.Range("K3").Value = 1212.8125
$ws.Range("L3").Value = 5830.3335
$ws.Range("M3").Value = -1098.8125
$ws.Range("N3").Value = -6058.3335

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H102").Value = 24999
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 24999
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 24999
$ws.Range("N102").Value = -31489
$ws.Range("M102").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1287.5946
$ws.Range("I58").Value = 1429.2632
$ws.Range("J58").Value = 1138.0555
$ws.Range("K58").Value = 1429.2632
$ws.Range("L58").Value = 1138.0555
$ws.Range("M58").Value = -1226.2632
$ws.Range("N58").Value = -1544.0555

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1738.5
$ws.Range("I99").Value = 1608.8572
$ws.Range("J99").Value = 1920
$ws.Range("K99").Value = 1608.8572
$ws.Range("L99").Value = 1920
$ws.Range("M99").Value = -110.8571999999999
$ws.Range("N99").Value = -4916

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 1738.5
$ws.Range("I126").Value = 1608.8572
$ws.Range("J126").Value = 1920
$ws.Range("K126").Value = 4826.571599999999
$ws.Range("L126").Value = 5760
$ws.Range("M126").Value = -2356.571599999999
$ws.Range("N126").Value = -10700

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1287.5946
$ws.Range("I136").Value = 1429.2632
$ws.Range("J136").Value = 1138.0555
$ws.Range("K136").Value = 4287.7896
$ws.Range("L136").Value = 3414.1665
$ws.Range("M136").Value = -1737.7896
$ws.Range("N136").Value = -8514.166499999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 11455
$ws.Range("I3").Value = 9576.666999999999
$ws.Range("J3").Value = 13333.333
$ws.Range("K3").Value = 28730.001
$ws.Range("L3").Value = 39999.999
$ws.Range("M3").Value = -28618.001
$ws.Range("N3").Value = -40223.999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 53750
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 53750
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 161250
$ws.Range("N17").Value = -161588
$ws.Range("M17").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 1291.9231
$ws.Range("I18").Value = 1117.7273
$ws.Range("J18").Value = 2250
$ws.Range("K18").Value = 3353.1819
$ws.Range("L18").Value = 6750
$ws.Range("M18").Value = -3184.1819
$ws.Range("N18").Value = -7088

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1960
$ws.Range("J34").Value = 1960
$ws.Range("L34").Value = 5880
$ws.Range("N34").Value = -6048

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 2156.8572
$ws.Range("J39").Value = 2156.8572
$ws.Range("L39").Value = 6470.571599999999
$ws.Range("N39").Value = -7058.571599999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 86133.664
$ws.Range("J55").Value = 93927.63
$ws.Range("L55").Value = 281782.89
$ws.Range("N55").Value = -282136.89

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 579.5599999999999
$ws.Range("I107").Value = 327.6
$ws.Range("J107").Value = 747.5333000000001
$ws.Range("K107").Value = 982.8000000000001
$ws.Range("L107").Value = 2242.5999
$ws.Range("M107").Value = 937.1999999999999
$ws.Range("N107").Value = -6082.5999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H130").Value = 1769.2307
$ws.Range("I130").Value = 500
$ws.Range("K130").Value = 1500
$ws.Range("M130").Value = 3520

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 849.0964
$ws.Range("I131").Value = 404.875
$ws.Range("J131").Value = 896.48
$ws.Range("K131").Value = 1214.625
$ws.Range("L131").Value = 2689.44
$ws.Range("M131").Value = 3825.375
$ws.Range("N131").Value = -12769.44

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 917.7222
$ws.Range("I132").Value = 465.64285
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 4190.78565
$ws.Range("L132").Value = 22500
$ws.Range("M132").Value = -1660.78565
$ws.Range("N132").Value = -27560

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1555.1052
$ws.Range("I113").Value = 1538.2727
$ws.Range("J113").Value = 1578.25
$ws.Range("K113").Value = 1538.2727
$ws.Range("L113").Value = 1578.25
$ws.Range("M113").Value = 631.7273
$ws.Range("N113").Value = -5918.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2000
$ws.Range("I40").Value = 2000
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 2000
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -1864
$ws.Range("N40").ClearContents()
